$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.573811
$ws.Range("H2").Value = 13.721433
$ws.Range("I2").Value = 0.1659009079913533
$ws.Range("J2").Value = 0.1659009079913533
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.03508633333333333
$ws.Range("N2").Value = 0.105259
$ws.Range("O2").Value = 0.004489537393262644
$ws.Range("P2").Value = 0.004489537393262644
$ws.Range("Q2").Value = 0.1604782573496667
$ws.Range("R2").Value = 1.444304316147
$ws.Range("S2").Value = 0.0007448183300034061
$ws.Range("T2").Value = 0.0007448183300034061

$ws.Range("G3").Value = 4.573811
$ws.Range("H3").Value = 13.721433
$ws.Range("I3").Value = 0.1659009079913533
$ws.Range("J3").Value = 0.1659009079913533
$ws.Range("M3").Value = 4.911922333333334
$ws.Range("O3").Value = 0.6285142074777995
$ws.Range("P3").Value = 0.6285142074777995
$ws.Range("Q3").Value = 22.46620439934567
$ws.Range("R3").Value = 202.195839594111
$ws.Range("S3").Value = 0.1042710777060328
$ws.Range("T3").Value = 0.1042710777060328

$ws.Range("G4").Value = 4.573811
$ws.Range("H4").Value = 13.721433
$ws.Range("I4").Value = 0.1659009079913533
$ws.Range("J4").Value = 0.1659009079913533
$ws.Range("M4").Value = 2.868124666666667
$ws.Range("N4").Value = 8.604374
$ws.Range("O4").Value = 0.3669962551289379
$ws.Range("P4").Value = 0.3669962551289379
$ws.Range("Q4").Value = 13.11826014977133
$ws.Range("R4").Value = 118.064341347942
$ws.Range("S4").Value = 0.06088501195531715
$ws.Range("T4").Value = 0.06088501195531715

$ws.Range("I5").Value = 0.5322852674812913
$ws.Range("J5").Value = 0.5322852674812913
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.03508633333333333
$ws.Range("N5").Value = 0.105259
$ws.Range("O5").Value = 0.004489537393262644
$ws.Range("P5").Value = 0.004489537393262644
$ws.Range("Q5").Value = 0.5148869477118889
$ws.Range("R5").Value = 4.633982529407
$ws.Range("S5").Value = 0.002389714612240066
$ws.Range("T5").Value = 0.002389714612240066

$ws.Range("I6").Value = 0.5322852674812913
$ws.Range("J6").Value = 0.5322852674812913
$ws.Range("M6").Value = 4.911922333333334
$ws.Range("O6").Value = 0.6285142074777995
$ws.Range("P6").Value = 0.6285142074777995
$ws.Range("Q6").Value = 72.0817611113879
$ws.Range("R6").Value = 648.7358500024911
$ws.Range("S6").Value = 0.3345488530431124
$ws.Range("T6").Value = 0.3345488530431124

$ws.Range("I7").Value = 0.5322852674812913
$ws.Range("J7").Value = 0.5322852674812913
$ws.Range("M7").Value = 2.868124666666667
$ws.Range("N7").Value = 8.604374
$ws.Range("O7").Value = 0.3669962551289379
$ws.Range("P7").Value = 0.3669962551289379
$ws.Range("Q7").Value = 42.08932125358911
$ws.Range("R7").Value = 378.803891282302
$ws.Range("S7").Value = 0.1953466998259389
$ws.Range("T7").Value = 0.1953466998259389

$ws.Range("G8").Value = 8.320867
$ws.Range("H8").Value = 24.962601
$ws.Range("I8").Value = 0.3018138245273554
$ws.Range("J8").Value = 0.3018138245273554
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.03508633333333333
$ws.Range("N8").Value = 0.105259
$ws.Range("O8").Value = 0.004489537393262644
$ws.Range("P8").Value = 0.004489537393262644
$ws.Range("Q8").Value = 0.2919487131843333
$ws.Range("R8").Value = 2.627538418659
$ws.Range("S8").Value = 0.001355004451019172
$ws.Range("T8").Value = 0.001355004451019172

$ws.Range("G9").Value = 8.320867
$ws.Range("H9").Value = 24.962601
$ws.Range("I9").Value = 0.3018138245273554
$ws.Range("J9").Value = 0.3018138245273554
$ws.Range("M9").Value = 4.911922333333334
$ws.Range("O9").Value = 0.6285142074777995
$ws.Range("P9").Value = 0.6285142074777995
$ws.Range("Q9").Value = 40.87145244999634
$ws.Range("R9").Value = 367.843072049967
$ws.Range("S9").Value = 0.1896942767286544
$ws.Range("T9").Value = 0.1896942767286544

$ws.Range("G10").Value = 8.320867
$ws.Range("H10").Value = 24.962601
$ws.Range("I10").Value = 0.3018138245273554
$ws.Range("J10").Value = 0.3018138245273554
$ws.Range("M10").Value = 2.868124666666667
$ws.Range("N10").Value = 8.604374
$ws.Range("O10").Value = 0.3669962551289379
$ws.Range("P10").Value = 0.3669962551289379
$ws.Range("Q10").Value = 23.86528389075266
$ws.Range("R10").Value = 214.787555016774
$ws.Range("S10").Value = 0.1107645433476818
$ws.Range("T10").Value = 0.1107645433476818
